$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 350
$ws.Range("F6").Value = 748
$ws.Range("F7").Value = 194
$ws.Range("F8").Value = 250
$ws.Range("F9").Value = 1048
$ws.Range("F11").Value = 307
$ws.Range("F12").Value = 605
$ws.Range("F13").Value = 168
$ws.Range("F24").Value = 201
$ws.Range("F25").Value = 15
$ws.Range("F26").Value = 147
$ws.Range("F27").Value = 581
$ws.Range("F28").Value = 954
$ws.Range("F30").Value = 205
$ws.Range("F31").Value = 1012
$ws.Range("F33").Value = 41
$ws.Range("F34").Value = 270

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 1026
$ws.Range("F5").Value = 1026
$ws.Range("F8").Value = 221
$ws.Range("F10").Value = 310
$ws.Range("F14").Value = 568
$ws.Range("F24").Value = 282
$ws.Range("F26").Value = 2236
$ws.Range("F31").Value = 8
$ws.Range("F33").Value = 95

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 2387
$ws.Range("F6").Value = 977
$ws.Range("F9").Value = 1227
$ws.Range("F10").Value = 322

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 2387
$ws.Range("F8").Value = 977
$ws.Range("F9").Value = 1227
$ws.Range("F10").Value = 322
$ws.Range("F13").Value = 350
$ws.Range("F14").Value = 748
$ws.Range("F15").Value = 194
$ws.Range("F17").Value = 250
$ws.Range("F18").Value = 1048
$ws.Range("F19").Value = 307
$ws.Range("F20").Value = 605
$ws.Range("F21").Value = 1026
$ws.Range("F31").Value = 201
$ws.Range("F32").Value = 147
$ws.Range("F33").Value = 581
$ws.Range("F34").Value = 954
$ws.Range("F35").Value = 568
$ws.Range("F38").Value = 205
$ws.Range("F43").Value = 282
$ws.Range("F44").Value = 282
$ws.Range("F46").Value = 1012
$ws.Range("F49").Value = 41
$ws.Range("F50").Value = 270
